# Add SMI dimension, update a few code table extra columns
$wb = $excel.ActiveWorkbook

# Remember the sheet that is active before we start, so we can restore it
# (selecting cells on the two sheets below will activate them in turn).
$origActive = $wb.ActiveSheet
$origActiveName = $origActive.Name

# ---------------------------------------------------------------------
# BehavioralHealthEvaluationType: add a 3rd column "BehavioralHealthEvaluationTypeCode"
# that mirrors the description column, except the "None" row (99998) which
# becomes "No Diagnosis" in the description column while the new code
# column keeps "None".
# ---------------------------------------------------------------------
$ws30 = $wb.Worksheets.Item("BehavioralHealthEvaluationType")

$ws30.Range("C1").Value = "BehavioralHealthEvaluationTypeCode"

for ($r = 2; $r -le 11; $r++) {
    $ws30.Cells.Item($r, 3).Value = $ws30.Cells.Item($r, 2).Value2
}

# Row 12 (99998): description changes from "None" to "No Diagnosis";
# the new code column keeps the original "None" value.
$ws30.Range("C12").Value = $ws30.Range("B12").Value2
$ws30.Range("B12").Value = "No Diagnosis"

# Row 13 (99999): "Unknown" carries through unchanged.
$ws30.Range("C13").Value = $ws30.Range("B13").Value2

$ws30.Columns.Item(3).ColumnWidth = 74.16666666666667

$null = $ws30.Activate()
$null = $ws30.Range("B13").Select()

# ---------------------------------------------------------------------
# MedicationType: add two columns, "MedicationTypeCategory" and
# "MedicationTypeCode", both mirroring the description column as-is.
# ---------------------------------------------------------------------
$ws31 = $wb.Worksheets.Item("MedicationType")

$ws31.Range("C1").Value = "MedicationTypeCategory"
$ws31.Range("D1").Value = "MedicationTypeCode"

for ($r = 2; $r -le 13; $r++) {
    $desc = $ws31.Cells.Item($r, 2).Value2
    $ws31.Cells.Item($r, 3).Value = $desc
    $ws31.Cells.Item($r, 4).Value = $desc
}

$ws31.Columns.Item(3).ColumnWidth = 24
$ws31.Columns.Item(4).ColumnWidth = 34.833333333333336

$null = $ws31.Activate()
$null = $ws31.Range("C14").Select()

# Restore the originally active sheet/tab.
$origSheet = $wb.Worksheets.Item($origActiveName)
$null = $origSheet.Activate()
